$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the data rows (2-29) in sorted-by-country order, inserting the
# new providers (Consortium GARR/IT, DigitalOcean Spaces/DE, Djezzy
# Alger+Annaba/DZ, Green Data/FR, HessFR/FR, Google Cloud Storage/WW).
$data = New-Object 'object[,]' 28,5
$data[0,0] = 'AF'
$data[0,1] = 'Afghan Wireless'
$data[0,2] = '10Gbps'
$data[0,3] = 'Afghan Wireless_10Gbps_AF.xml'
$data[0,4] = 'No'
$data[1,0] = 'AR'
$data[1,1] = 'ClearNet S.A.'
$data[1,2] = '10Gbps'
$data[1,3] = 'ClearNet S.A._10Gbps_AR.xml'
$data[1,4] = 'No'
$data[2,0] = 'CH'
$data[2,1] = '31173 Service - Zurich'
$data[2,2] = '10Gbps'
$data[2,3] = '31173 Services AB_10Gbps_Zurich_CH'
$data[2,4] = 'No'
$data[3,0] = 'DE'
$data[3,1] = '31173 Service - Frankfurt'
$data[3,2] = '10Gbps'
$data[3,3] = '31173 Services AB_10Gbps_Frankfurt_DE.xml'
$data[3,4] = 'No'
$data[4,0] = 'DE'
$data[4,1] = 'DigitalOcean Spaces'
$data[4,2] = 'CDN'
$data[4,3] = 'DigitalOcean Spaces_Germany_Frankfurt_DE.xml'
$data[4,4] = 'No'
$data[5,0] = 'DZ'
$data[5,1] = 'Djezzy Alger'
$data[5,2] = '1Gbps'
$data[5,3] = 'Djezzy Alger_1Gbps_DZ.xml'
$data[5,4] = 'No'
$data[6,0] = 'DZ'
$data[6,1] = 'Djezzy Annaba'
$data[6,2] = '1Gbps'
$data[6,3] = 'Djezzy Annaba_1Gbps_DZ.xml'
$data[6,4] = 'No'
$data[7,0] = 'FR'
$data[7,1] = 'Appliwave BBR'
$data[7,2] = '40Gbps'
$data[7,3] = 'Appliwave_BBR_40Gbps_FR.xml'
$data[7,4] = 'No'
$data[8,0] = 'FR'
$data[8,1] = 'Bouygues Telecom Paris'
$data[8,2] = '20Gbps'
$data[8,3] = 'Bouygues Telecom FR_20Gbps_FR.xml'
$data[8,4] = 'No'
$data[9,0] = 'FR'
$data[9,1] = 'Bouygues Telecom Lille'
$data[9,2] = '10Gbps'
$data[9,3] = 'Bouygues Telecom Lille_10Gbps_FR.xml'
$data[9,4] = 'No'
$data[10,0] = 'FR'
$data[10,1] = 'Bouygues Telecom Lyon'
$data[10,2] = '10Gbps'
$data[10,3] = 'Bouygues Telecom Lyon_10Gbps_FR.xml'
$data[10,4] = 'No'
$data[11,0] = 'FR'
$data[11,1] = 'Bouygues Telecom Marseille BBR'
$data[11,2] = '10Gbps'
$data[11,3] = 'Bouygues_Telecom_MRS_BBR_10Gbps_FR.xml'
$data[11,4] = 'No'
$data[12,0] = 'FR'
$data[12,1] = 'Bouygues Telecom Marseille'
$data[12,2] = '10Gbps'
$data[12,3] = 'Bouygues Telecom Marseille_10Gbps_FR.xml'
$data[12,4] = 'No'
$data[13,0] = 'FR'
$data[13,1] = 'Bouygues Telecom Poitiers'
$data[13,2] = '10Gbps'
$data[13,3] = 'Bouygues Telecom Poitiers_10Gbps_FR.xml'
$data[13,4] = 'No'
$data[14,0] = 'FR'
$data[14,1] = 'Bouygues Telecom Rennes'
$data[14,2] = '10Gbps'
$data[14,3] = 'Bouygues Telecom Rennes_10Gbps_FR.xml'
$data[14,4] = 'No'
$data[15,0] = 'FR'
$data[15,1] = 'Bouygues Telecom Strasbourg'
$data[15,2] = '10Gbps'
$data[15,3] = 'Bouygues Telecom Strasbourg_10Gbps_FR.xml'
$data[15,4] = 'No'
$data[16,0] = 'FR'
$data[16,1] = 'Bouygue Telecom Toulouse'
$data[16,2] = '10Gbps'
$data[16,3] = 'Bouygues Telecom_10Gbps_FR.xml'
$data[16,4] = 'No'
$data[17,0] = 'FR'
$data[17,1] = 'Green Data'
$data[17,2] = '10Gbps'
$data[17,3] = 'Green Data_10Gbps_FR.xml'
$data[17,4] = 'No'
$data[18,0] = 'FR'
$data[18,1] = 'HessFR'
$data[18,2] = '10Gbps'
$data[18,3] = 'HessFR_10Gbps_FR.xml'
$data[18,4] = 'No'
$data[19,0] = 'GB'
$data[19,1] = '31173 Service - London'
$data[19,2] = '10Gbps'
$data[19,3] = '31173 Services AB_10Gbps_London_GB.xml'
$data[19,4] = 'No'
$data[20,0] = 'GB'
$data[20,1] = 'Amazon S3 - Eu-west2'
$data[20,2] = 'CDN'
$data[20,3] = 'Academic Computer Club_20Gbps_Umea_SE.xml'
$data[20,4] = 'No'
$data[21,0] = 'IT'
$data[21,1] = 'Consortium GARR'
$data[21,2] = '100Gbps'
$data[21,3] = 'Consortium GARR_100Gbps_IT.xml'
$data[21,4] = 'No'
$data[22,0] = 'NL'
$data[22,1] = '31173 Service - Amsterdam'
$data[22,2] = '10Gbps'
$data[22,3] = '31173 Services AB_10Gbps_Amsterdam_NL.xml'
$data[22,4] = 'No'
$data[23,0] = 'SE'
$data[23,1] = 'Bahnhof AB'
$data[23,2] = '10Gbps'
$data[23,3] = 'Bahnhof AB_10Gbps_SE.xml'
$data[23,4] = 'No'
$data[24,0] = 'SE'
$data[24,1] = 'Academic Computer Club - Umea'
$data[24,2] = '20Gbps'
$data[24,3] = 'Academic Computer Club_20Gbps_Umea_SE.xml'
$data[24,4] = 'No'
$data[25,0] = 'US'
$data[25,1] = 'Amazon S3 - US-West1'
$data[25,2] = 'CDN'
$data[25,3] = 'AmazonS3_US-West-1_US.xml'
$data[25,4] = 'No'
$data[26,0] = 'WW'
$data[26,1] = 'Cloudflare Blender'
$data[26,2] = 'CDN'
$data[26,3] = 'Cloudflare_Blender_CDN.xml'
$data[26,4] = 'No'
$data[27,0] = 'WW'
$data[27,1] = 'Google Cloud Storage'
$data[27,2] = 'CDN'
$data[27,3] = 'Google Cloud Storage_CDN.xml'
$data[27,4] = 'No'
$ws.Range("A2:E29").Value = $data

# New used range / dimension grows to A1:E29
$ws.Range("A1:E29").Select() | Out-Null

# Conditional formatting on column E (whole column) highlighting Yes/No,
# matching Excel's built-in "Highlight Cells Rules > Equal To" red/green
# styles.
$cfRange = $ws.Range("E1:E1048576")

$fcNo = $cfRange.FormatConditions.Add(8, 3, '"No"')
$fcNo.Font.Color = 393372
$fcNo.Interior.Color = 13551615

$fcYes = $cfRange.FormatConditions.Add(8, 3, '"Yes"')
$fcYes.Font.Color = 24832
$fcYes.Interior.Color = 13561798

# Restore the selection Excel shows after the edit.
$ws.Range("D21").Select() | Out-Null
